# Automatische test-sync: 2025-06-22 18:57:50
#
# Applies:
#   1. Dashboard sheet: re-order several category rows (A3:A12) and append a
#      new row 15 ("Juridisch / Contract" = 1).
#   2. Chart1 (embedded on Dashboard): extend the category/value series
#      references from row 14 to row 15.
#   3. Logs sheet: append a new log row (row 27) for the GDPR / privacy
#      question, and extend the conditional-formatting ranges on columns D
#      and G to cover the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Dashboard: re-order categories (values in column B stay put, only the
#    labels in column A are shuffled) and add the new "Juridisch / Contract"
#    row.
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value  = "Afmelding / Nieuwsbrief"
$dash.Range("A4").Value  = "Productinformatie"
$dash.Range("A5").Value  = "Samenwerking / Partnerverzoek"
$dash.Range("A6").Value  = "Offerte / Prijsaanvraag"
$dash.Range("A7").Value  = "Sollicitatie / Vacature"
$dash.Range("A9").Value  = "Retour / Terugbetaling"
$dash.Range("A11").Value = "Uitnodiging / Evenement"
$dash.Range("A12").Value = "Openingstijden / Locatie"

$dash.Range("A15").Value = "Juridisch / Contract"
$dash.Range("B15").Value = 1

# ---------------------------------------------------------------------------
# 2. Chart series: stretch the category/value ranges to row 15.
# ---------------------------------------------------------------------------
$chart = $dash.ChartObjects(1).Chart
$ser = $chart.SeriesCollection(1)
$ser.XValues = "'Dashboard'!`$A`$2:`$A`$15"
$ser.Values  = "'Dashboard'!`$B`$2:`$B`$15"

# ---------------------------------------------------------------------------
# 3. Logs: append the new row 27.
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A27").Value = "Privacybeleid"
$logs.Range("B27").Value = "mailmind.test@zohomail.eu"
$logs.Range("C27").Value = "Wat doen jullie met klantgegevens volgens GDPR?"
$logs.Range("D27").Value = "Juridisch / Contract"

$reply = @'
Beste klant,
Dank u voor uw vraag over hoe wij omgaan met klantgegevens volgens de GDPR. Wij nemen de privacy en bescherming van persoonlijke gegevens van onze klanten uiterst serieus en volgen hierbij strikt de richtlijnen van de Algemene Verordening Gegevensbescherming (AVG/GDPR).
Om uw privacy te waarborgen, zorgen wij ervoor dat alle klantgegevens veilig worden opgeslagen en verwerkt volgens de geldende wet- en regelgeving. Dit houdt in dat we alleen persoonlijke gegevens verzamelen die relevant zijn voor het leveren van onze diensten en dat we deze gegevens niet delen met derden zonder uw toestemming, tenzij dit noodzakelijk is voor de uitvoering van onze diensten.
Mocht u nog specifieke vragen hebben over hoe wij omgaan met uw gegevens of wilt u meer informatie ontvangen, aarzel dan niet om contact met ons op te nemen.
Met vriendelijke groet,
[Naam] Nederlandse e-mailassistent 
[Bedrijfsnaam]
'@
$logs.Range("E27").Value = $reply

$logs.Range("F27").Value = "2025-06-22 18:57:33"
$logs.Range("G27").Value = "Ja"

# Writing the multi-paragraph reply above flags the row for an explicit
# "ht"/"customHeight" override (simulating Excel's content re-measure).
# None of the other log rows carry an explicit height, so re-run autofit to
# drop back to the sheet's standard (implicit) row height.
$logs.Rows.Item(27).AutoFit()

# ---------------------------------------------------------------------------
# Extend the conditional formatting ranges on Logs (D2:D26 -> D2:D27 and
# G2:G26 -> G2:G27) so the new row gets the same colour rules. Priorities /
# dxfIds are preserved because we modify the existing rules in place instead
# of deleting and re-creating them.
# ---------------------------------------------------------------------------
$dFC = $logs.Range("D2:D26").FormatConditions
for ($i = 1; $i -le $dFC.Count; $i++) {
    $dFC.Item($i).ModifyAppliesToRange($logs.Range("D2:D27"))
}

$gFC = $logs.Range("G2:G26").FormatConditions
for ($i = 1; $i -le $gFC.Count; $i++) {
    $gFC.Item($i).ModifyAppliesToRange($logs.Range("G2:G27"))
}
